$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 493-527 (columns D,H,I,J,K,L,M,O,P) ---
$ws.Range("D493").Value = 44826
$ws.Range("H493").Value = "Copenhague"
$ws.Range("I493").Value = "Primera"
$ws.Range("J493").Value = 250
$ws.Range("K493").Value = 2500
$ws.Range("L493").Value = 2500
$ws.Range("M493").Value = 2500
$ws.Range("O493").Value = "Región Metropolitana"
$ws.Range("P493").Value = 2500

$ws.Range("D494").Value = 44826
$ws.Range("H494").Value = "Crespo record"
$ws.Range("I494").Value = "Primera"
$ws.Range("J494").Value = 250
$ws.Range("K494").Value = 2200
$ws.Range("L494").Value = 2200
$ws.Range("M494").Value = 2200
$ws.Range("O494").Value = "Región Metropolitana"
$ws.Range("P494").Value = 2200

$ws.Range("D495").Value = 44413
$ws.Range("H495").Value = "Copenhague"
$ws.Range("I495").Value = "Primera"
$ws.Range("J495").Value = 250
$ws.Range("K495").Value = 1400
$ws.Range("L495").Value = 1400
$ws.Range("M495").Value = 1400
$ws.Range("O495").Value = "Región Metropolitana"
$ws.Range("P495").Value = 1400

$ws.Range("D496").Value = 44413
$ws.Range("H496").Value = "Crespo record"
$ws.Range("I496").Value = "Primera"
$ws.Range("J496").Value = 250
$ws.Range("K496").Value = 1200
$ws.Range("L496").Value = 1200
$ws.Range("M496").Value = 1200
$ws.Range("O496").Value = "Región Metropolitana"
$ws.Range("P496").Value = 1200

$ws.Range("D497").Value = 44515
$ws.Range("H497").Value = "Crespo record"
$ws.Range("I497").Value = "Primera"
$ws.Range("J497").Value = 250
$ws.Range("K497").Value = 1200
$ws.Range("L497").Value = 1200
$ws.Range("M497").Value = 1200
$ws.Range("O497").Value = "Región Metropolitana"
$ws.Range("P497").Value = 1200

$ws.Range("D498").Value = 44515
$ws.Range("H498").Value = "Crespo record"
$ws.Range("I498").Value = "Segunda"
$ws.Range("J498").Value = 250
$ws.Range("K498").Value = 1000
$ws.Range("L498").Value = 1000
$ws.Range("M498").Value = 1000
$ws.Range("O498").Value = "Región Metropolitana"
$ws.Range("P498").Value = 1000

$ws.Range("D499").Value = 44356
$ws.Range("H499").Value = "Copenhague"
$ws.Range("I499").Value = "Primera"
$ws.Range("J499").Value = 100
$ws.Range("K499").Value = 1200
$ws.Range("L499").Value = 1200
$ws.Range("M499").Value = 1200
$ws.Range("O499").Value = "Región Metropolitana"
$ws.Range("P499").Value = 1200

$ws.Range("D500").Value = 44379
$ws.Range("H500").Value = "Copenhague"
$ws.Range("I500").Value = "Primera"
$ws.Range("J500").Value = 600
$ws.Range("K500").Value = 1200
$ws.Range("L500").Value = 1200
$ws.Range("M500").Value = 1200
$ws.Range("O500").Value = "Región Metropolitana"
$ws.Range("P500").Value = 1200

$ws.Range("D501").Value = 44379
$ws.Range("H501").Value = "Crespo record"
$ws.Range("I501").Value = "Primera"
$ws.Range("J501").Value = 800
$ws.Range("K501").Value = 1000
$ws.Range("L501").Value = 1200
$ws.Range("M501").Value = 1100
$ws.Range("O501").Value = "Región Metropolitana"
$ws.Range("P501").Value = 1100

$ws.Range("D502").Value = 44322
$ws.Range("H502").Value = "Crespo record"
$ws.Range("I502").Value = "Primera"
$ws.Range("J502").Value = 250
$ws.Range("K502").Value = 1200
$ws.Range("L502").Value = 1200
$ws.Range("M502").Value = 1200
$ws.Range("O502").Value = "Región Metropolitana"
$ws.Range("P502").Value = 1200

$ws.Range("D503").Value = 44322
$ws.Range("H503").Value = "Crespo record"
$ws.Range("I503").Value = "Segunda"
$ws.Range("J503").Value = 250
$ws.Range("K503").Value = 1000
$ws.Range("L503").Value = 1000
$ws.Range("M503").Value = 1000
$ws.Range("O503").Value = "Región Metropolitana"
$ws.Range("P503").Value = 1000

$ws.Range("D504").Value = 44497
$ws.Range("H504").Value = "Crespo record"
$ws.Range("I504").Value = "Primera"
$ws.Range("J504").Value = 750
$ws.Range("K504").Value = 1200
$ws.Range("L504").Value = 1200
$ws.Range("M504").Value = 1200
$ws.Range("O504").Value = "Región de Coquimbo"
$ws.Range("P504").Value = 1200

$ws.Range("D505").Value = 44782
$ws.Range("H505").Value = "Crespo record"
$ws.Range("I505").Value = "Primera"
$ws.Range("J505").Value = 1400
$ws.Range("K505").Value = 2000
$ws.Range("L505").Value = 2000
$ws.Range("M505").Value = 2000
$ws.Range("O505").Value = "Región Metropolitana"
$ws.Range("P505").Value = 2000

$ws.Range("D506").Value = 44435
$ws.Range("H506").Value = "Copenhague"
$ws.Range("I506").Value = "Primera"
$ws.Range("J506").Value = 1100
$ws.Range("K506").Value = 1400
$ws.Range("L506").Value = 1400
$ws.Range("M506").Value = 1400
$ws.Range("O506").Value = "Región Metropolitana"
$ws.Range("P506").Value = 1400

$ws.Range("D507").Value = 44435
$ws.Range("H507").Value = "Crespo record"
$ws.Range("I507").Value = "Primera"
$ws.Range("J507").Value = 1700
$ws.Range("K507").Value = 1200
$ws.Range("L507").Value = 1400
$ws.Range("M507").Value = 1288
$ws.Range("O507").Value = "Región Metropolitana"
$ws.Range("P507").Value = 1288

$ws.Range("D508").Value = 44435
$ws.Range("H508").Value = "Crespo record"
$ws.Range("I508").Value = "Segunda"
$ws.Range("J508").Value = 1200
$ws.Range("K508").Value = 1000
$ws.Range("L508").Value = 1000
$ws.Range("M508").Value = 1000
$ws.Range("O508").Value = "Región Metropolitana"
$ws.Range("P508").Value = 1000

$ws.Range("D509").Value = 44435
$ws.Range("H509").Value = "Crespo record"
$ws.Range("I509").Value = "Segunda"
$ws.Range("J509").Value = 250
$ws.Range("K509").Value = 1000
$ws.Range("L509").Value = 1000
$ws.Range("M509").Value = 1000
$ws.Range("O509").Value = "Región del Maule"
$ws.Range("P509").Value = 1000

$ws.Range("D510").Value = 44251
$ws.Range("H510").Value = "Crespo record"
$ws.Range("I510").Value = "Primera"
$ws.Range("J510").Value = 250
$ws.Range("K510").Value = 1500
$ws.Range("L510").Value = 1500
$ws.Range("M510").Value = 1500
$ws.Range("O510").Value = "Provincia de Quillota"
$ws.Range("P510").Value = 1500

$ws.Range("D511").Value = 44319
$ws.Range("H511").Value = "Crespo record"
$ws.Range("I511").Value = "Segunda"
$ws.Range("J511").Value = 250
$ws.Range("K511").Value = 1000
$ws.Range("L511").Value = 1000
$ws.Range("M511").Value = 1000
$ws.Range("O511").Value = "Región Metropolitana"
$ws.Range("P511").Value = 1000

$ws.Range("D512").Value = 44344
$ws.Range("H512").Value = "Crespo record"
$ws.Range("I512").Value = "Primera"
$ws.Range("J512").Value = 1400
$ws.Range("K512").Value = 1000
$ws.Range("L512").Value = 1000
$ws.Range("M512").Value = 1000
$ws.Range("O512").Value = "Región Metropolitana"
$ws.Range("P512").Value = 1000

$ws.Range("D513").Value = 44232
$ws.Range("H513").Value = "Copenhague"
$ws.Range("I513").Value = "Primera"
$ws.Range("J513").Value = 700
$ws.Range("K513").Value = 1600
$ws.Range("L513").Value = 1600
$ws.Range("M513").Value = 1600
$ws.Range("O513").Value = "Región de Coquimbo"
$ws.Range("P513").Value = 1600

$ws.Range("D514").Value = 44232
$ws.Range("H514").Value = "Crespo record"
$ws.Range("I514").Value = "Primera"
$ws.Range("J514").Value = 700
$ws.Range("K514").Value = 1600
$ws.Range("L514").Value = 1600
$ws.Range("M514").Value = 1600
$ws.Range("O514").Value = "Región de Coquimbo"
$ws.Range("P514").Value = 1600

$ws.Range("D515").Value = 44504
$ws.Range("H515").Value = "Crespo record"
$ws.Range("I515").Value = "Primera"
$ws.Range("J515").Value = 500
$ws.Range("K515").Value = 1200
$ws.Range("L515").Value = 1200
$ws.Range("M515").Value = 1200
$ws.Range("O515").Value = "Región de Coquimbo"
$ws.Range("P515").Value = 1200

$ws.Range("D516").Value = 44484
$ws.Range("H516").Value = "Crespo record"
$ws.Range("I516").Value = "Primera"
$ws.Range("J516").Value = 1200
$ws.Range("K516").Value = 1200
$ws.Range("L516").Value = 1200
$ws.Range("M516").Value = 1200
$ws.Range("O516").Value = "Región Metropolitana"
$ws.Range("P516").Value = 1200

$ws.Range("D517").Value = 44665
$ws.Range("H517").Value = "Crespo record"
$ws.Range("I517").Value = "Primera"
$ws.Range("J517").Value = 500
$ws.Range("K517").Value = 1800
$ws.Range("L517").Value = 1800
$ws.Range("M517").Value = 1800
$ws.Range("O517").Value = "Región Metropolitana"
$ws.Range("P517").Value = 1800

$ws.Range("D518").Value = 44452
$ws.Range("H518").Value = "Crespo record"
$ws.Range("I518").Value = "Segunda"
$ws.Range("J518").Value = 500
$ws.Range("K518").Value = 1000
$ws.Range("L518").Value = 1000
$ws.Range("M518").Value = 1000
$ws.Range("O518").Value = "Región del Maule"
$ws.Range("P518").Value = 1000

$ws.Range("D519").Value = 44189
$ws.Range("H519").Value = "Copenhague"
$ws.Range("I519").Value = "Primera"
$ws.Range("J519").Value = 500
$ws.Range("K519").Value = 1400
$ws.Range("L519").Value = 1500
$ws.Range("M519").Value = 1450
$ws.Range("O519").Value = "Región Metropolitana"
$ws.Range("P519").Value = 1450

$ws.Range("D520").Value = 44189
$ws.Range("H520").Value = "Crespo record"
$ws.Range("I520").Value = "Primera"
$ws.Range("J520").Value = 500
$ws.Range("K520").Value = 1100
$ws.Range("L520").Value = 1200
$ws.Range("M520").Value = 1150
$ws.Range("O520").Value = "Región del Maule"
$ws.Range("P520").Value = 1150

$ws.Range("D521").Value = 44701
$ws.Range("H521").Value = "Crespo record"
$ws.Range("I521").Value = "Primera"
$ws.Range("J521").Value = 600
$ws.Range("K521").Value = 1800
$ws.Range("L521").Value = 1800
$ws.Range("M521").Value = 1800
$ws.Range("O521").Value = "Región Metropolitana"
$ws.Range("P521").Value = 1800

$ws.Range("D522").Value = 44701
$ws.Range("H522").Value = "Crespo record"
$ws.Range("I522").Value = "Segunda"
$ws.Range("J522").Value = 600
$ws.Range("K522").Value = 1600
$ws.Range("L522").Value = 1600
$ws.Range("M522").Value = 1600
$ws.Range("O522").Value = "Región Metropolitana"
$ws.Range("P522").Value = 1600

$ws.Range("D523").Value = 44516
$ws.Range("H523").Value = "Copenhague"
$ws.Range("I523").Value = "Primera"
$ws.Range("J523").Value = 600
$ws.Range("K523").Value = 1400
$ws.Range("L523").Value = 1400
$ws.Range("M523").Value = 1400
$ws.Range("O523").Value = "Región Metropolitana"
$ws.Range("P523").Value = 1400

$ws.Range("D524").Value = 44516
$ws.Range("H524").Value = "Crespo record"
$ws.Range("I524").Value = "Primera"
$ws.Range("J524").Value = 500
$ws.Range("K524").Value = 1200
$ws.Range("L524").Value = 1200
$ws.Range("M524").Value = 1200
$ws.Range("O524").Value = "Región Metropolitana"
$ws.Range("P524").Value = 1200

$ws.Range("D525").Value = 44516
$ws.Range("H525").Value = "Crespo record"
$ws.Range("I525").Value = "Segunda"
$ws.Range("J525").Value = 500
$ws.Range("K525").Value = 1000
$ws.Range("L525").Value = 1000
$ws.Range("M525").Value = 1000
$ws.Range("O525").Value = "Región Metropolitana"
$ws.Range("P525").Value = 1000

$ws.Range("D526").Value = 44186
$ws.Range("H526").Value = "Crespo record"
$ws.Range("I526").Value = "Primera"
$ws.Range("J526").Value = 300
$ws.Range("K526").Value = 1000
$ws.Range("L526").Value = 1200
$ws.Range("M526").Value = 1100
$ws.Range("O526").Value = "Región del Maule"
$ws.Range("P526").Value = 1100

$ws.Range("D527").Value = 44463
$ws.Range("H527").Value = "Copenhague"
$ws.Range("I527").Value = "Primera"
$ws.Range("J527").Value = 600
$ws.Range("K527").Value = 1500
$ws.Range("L527").Value = 1500
$ws.Range("M527").Value = 1500
$ws.Range("O527").Value = "Región Metropolitana"
$ws.Range("P527").Value = 1500

# --- Insert two new rows at 528-529 (old row 528 shifts to 530) ---
$ws.Range("528:529").EntireRow.Insert()

# --- Populate new rows 528 and 529 ---
$ws.Range("A528").Value = 4
$ws.Range("B528").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C528").Value = "Los Lagos"
$ws.Range("D528").Value = 44463
$ws.Range("E528").Value = 10
$ws.Range("F528").Value = 100112006
$ws.Range("G528").Value = "Repollo"
$ws.Range("H528").Value = "Crespo record"
$ws.Range("I528").Value = "Primera"
$ws.Range("J528").Value = 500
$ws.Range("K528").Value = 1200
$ws.Range("L528").Value = 1200
$ws.Range("M528").Value = 1200
$ws.Range("N528").Value = "$/unidad"
$ws.Range("O528").Value = "Región Metropolitana"
$ws.Range("P528").Value = 1200
$ws.Range("Q528").Value = 1
$ws.Range("R528").Value = "Hortaliza"

$ws.Range("A529").Value = 4
$ws.Range("B529").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C529").Value = "Los Lagos"
$ws.Range("D529").Value = 44463
$ws.Range("E529").Value = 10
$ws.Range("F529").Value = 100112006
$ws.Range("G529").Value = "Repollo"
$ws.Range("H529").Value = "Crespo record"
$ws.Range("I529").Value = "Segunda"
$ws.Range("J529").Value = 500
$ws.Range("K529").Value = 1000
$ws.Range("L529").Value = 1000
$ws.Range("M529").Value = 1000
$ws.Range("N529").Value = "$/unidad"
$ws.Range("O529").Value = "Región del Maule"
$ws.Range("P529").Value = 1000
$ws.Range("Q529").Value = 1
$ws.Range("R529").Value = "Hortaliza"

